$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings stay as text (matches source formatting)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "43.027.02"
$ws.Range("E2").Value = "  +1.80%  "
$ws.Range("D3").Value = "2.307.70"
$ws.Range("E3").Value = "  +1.45%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "302.66"
$ws.Range("E5").Value = "  +0.90%  "
$ws.Range("D6").Value = "101.20"
$ws.Range("E6").Value = "  +4.94%  "
$ws.Range("D7").Value = "0.505"
$ws.Range("E7").Value = "  +1.68%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("E9").Value = "  +4.96%  "
$ws.Range("D10").Value = "35.38"
$ws.Range("E10").Value = "  +6.05%  "
$ws.Range("E11").Value = "  +0.88%  "
$ws.Range("E12").Value = "  +3.68%  "
$ws.Range("D13").Value = "18.00"
$ws.Range("E13").Value = "  +15.38%  "
$ws.Range("D14").Value = "6.92"
$ws.Range("E14").Value = "  +3.57%  "
$ws.Range("D15").Value = "2.683.99"
$ws.Range("E15").Value = "  +2.20%  "
$ws.Range("D16").Value = "2.283.85"
$ws.Range("E16").Value = "  -0.53%  "
$ws.Range("D17").Value = "0.813"
$ws.Range("E17").Value = "  +3.54%  "
$ws.Range("D18").Value = "42.927.19"
$ws.Range("E18").Value = "  +1.76%  "
$ws.Range("D19").Value = "12.62"
$ws.Range("E19").Value = "  +7.71%  "
$ws.Range("E20").Value = "  +3.02%  "
$ws.Range("D21").Value = "0.0₃0905"
$ws.Range("E21").Value = "  +1.56%  "
$ws.Range("D22").Value = "67.95"
$ws.Range("E22").Value = "  +2.14%  "
$ws.Range("D23").Value = "237.26"
$ws.Range("E23").Value = "  +0.73%  "
$ws.Range("E24").Value = "  +10.99%  "
$ws.Range("E25").Value = "  +0.07%  "
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  -0.13%  "
$ws.Range("D27").Value = "24.81"
$ws.Range("E27").Value = "  +3.26%  "
$ws.Range("D28").Value = "167.49"
$ws.Range("E28").Value = "  -0.82%  "
$ws.Range("E29").Value = "  +0.62%  "
$ws.Range("D30").Value = "34.14"
$ws.Range("E30").Value = "  +1.18%  "
$ws.Range("D31").Value = "9.25"
$ws.Range("E31").Value = "  +0.51%  "
$ws.Range("E32").Value = "  +0.11%  "
$ws.Range("D33").Value = "5.03"
$ws.Range("E33").Value = "  +2.26%  "
$ws.Range("D34").Value = "4.63"
$ws.Range("E34").Value = "  +2.59%  "
$ws.Range("E35").Value = "  +3.64%  "
$ws.Range("D36").Value = "16.98"
$ws.Range("E36").Value = "  +1.81%  "
$ws.Range("D37").Value = "0.0692"
$ws.Range("E37").Value = "  +0.60%  "
$ws.Range("E38").Value = "  +3.00%  "
$ws.Range("D39").Value = "2.83"
$ws.Range("E39").Value = "  +1.44%  "
$ws.Range("E40").Value = "  +3.91%  "
$ws.Range("E41").Value = "  +1.00%  "
$ws.Range("D42").Value = "2.31"
$ws.Range("E42").Value = "  -4.58%  "
$ws.Range("D43").Value = "2.002.59"
$ws.Range("E43").Value = "  +1.56%  "
$ws.Range("D44").Value = "0.0288"
$ws.Range("E44").Value = "  +3.70%  "
$ws.Range("D45").Value = "10.27"
$ws.Range("E45").Value = "  +7.29%  "
$ws.Range("D46").Value = "17.55"
$ws.Range("E46").Value = "  +1.09%  "
$ws.Range("D47").Value = "2.86"
$ws.Range("E47").Value = "  +2.53%  "
$ws.Range("D48").Value = "55.85"
$ws.Range("E48").Value = "  +6.71%  "
$ws.Range("D49").Value = "2.526.38"
$ws.Range("E49").Value = "  +1.14%  "
$ws.Range("E50").Value = "  +3.03%  "
$ws.Range("D51").Value = "4.56"
$ws.Range("E51").Value = "  +1.16%  "
